# Update marksheet corrected/total marks in the "quiz" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: B11 3 -> 5
$ws.Range("B11").Value = 5

# Total row: B12 51 -> 85
$ws.Range("B12").Value = 85

# Total row: E12 text "48/84" -> "85/140"
$ws.Range("E12").Value = "85/140"
